# Fix the "b9d2_unkown" typo -> "b9d2_unknown" across the genotype grid
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("genotype")
$ws.Range("B2:L9").Value = "b9d2_unknown"

# Switch the active/selected sheet from "temperature" to "genotype",
# and move the selection on genotype from D12 to B2
$ws.Activate()
$ws.Range("B2").Select()
